# Insert a new weekly price record for "Albahaca" (Feria Lagunitas de Puerto
# Montt) as row 47, pushing the existing rows 47-93 down to 48-94.
#
# Excel's Rows.Insert() shifts all data/styles/dimension automatically, so we
# insert the row, copy the (now shifted) old-row-47 contents from row 48 into
# the freshly inserted row 47 (to carry over the repeated columns such as
# Mercado, Region, Categoria, Unidad de comercializacion, Origen, etc.), and
# finally overwrite the handful of cells that hold the new record's unique
# values (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 47; rows 47:93 move to 48:94.
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the same repeated values as the row that used
# to be row 47 and now lives at row 48 (Mercado ID, Mercado, Region, Codreg,
# Categoria ID, Categoria, Variedad, Calidad, Unidad de comercializacion,
# Origen, Kg o Unidades, Clasificacion).
$ws.Range("A48:R48").Copy()
$ws.Range("A47").PasteSpecial()

# Overwrite the new record's own values.
$ws.Range("D47").Value = 44589
$ws.Range("J47").Value = 150
$ws.Range("K47").Value = 6000
$ws.Range("L47").Value = 6000
$ws.Range("M47").Value = 6000
$ws.Range("P47").Value = 1000
